$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Rows 15-16 swap Avalanche <-> WrappedliquidstakedEther2.0; rows 48-49 swap Mantle <-> ImmutableX.
# Column D values that look numeric are prefixed with a leading apostrophe so Excel stores them
# as text (matching the source data format, e.g. "3.220.12") instead of coercing to a float.

$ws.Range("D2").Value = "89.913.14"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "3.220.12"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  +6.00%  "
$ws.Range("D6").Value = "'629.40"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("E7").Value = "  +6.31%  "
$ws.Range("D8").Value = "'0.697"
$ws.Range("E8").Value = "  +5.50%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.216.81"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "'0.573"
$ws.Range("E11").Value = "  +6.27%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  +6.34%  "
$ws.Range("D14").Value = "'5.43"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'33.52"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.814.87"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "89.635.71"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("D18").Value = "3.200.40"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "'0.0000235"
$ws.Range("E19").Value = "  +80.36%  "
$ws.Range("D20").Value = "'3.47"
$ws.Range("E20").Value = "  +17.27%  "
$ws.Range("D21").Value = "'13.57"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "'439.97"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("D23").Value = "'8.65"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "'5.10"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'5.23"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "'11.95"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("D27").Value = "'81.49"
$ws.Range("E27").Value = "  +10.92%  "
$ws.Range("D28").Value = "3.390.53"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'0.160"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'4.15"
$ws.Range("E32").Value = "  +37.28%  "
$ws.Range("D33").Value = "'8.58"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "'544.61"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'7.02"
$ws.Range("E35").Value = "  +5.85%  "
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "'22.43"
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'146.67"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'173.79"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").Value = "'1.26"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.747"
$ws.Range("E49").Value = "  +8.14%  "
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("E51").Value = "  +5.64%  "
